$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.074.23"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "2.758.40"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'575.37"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "'158.99"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  -3.60%  "
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").Value = "'0.385"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("D12").Value = "'5.68"
$ws.Range("E12").Value = "  -16.23%  "
$ws.Range("D13").Value = "3.248.40"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "'26.88"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "63.688.31"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "'0.0000151"
$ws.Range("E16").Value = "  -2.75%  "
$ws.Range("D17").Value = "2.762.99"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "'12.15"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'4.86"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").Value = "'356.17"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  -4.17%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").Value = "'65.45"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").Value = "'8.60"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "0.0₃0908"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("E29").Value = "  -3.33%  "
$ws.Range("D30").Value = "'7.23"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").Value = "'1.24"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").Value = "'169.61"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").Value = "'20.27"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").Value = "'4.92"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").Value = "'1.48"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "'6.29"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'340.54"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'4.20"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").Value = "'39.21"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").Value = "'21.46"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").Value = "'21.73"
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("D45").Value = "'0.0589"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("D46").Value = "'0.634"
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("D47").Value = "'0.0255"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").Value = "'135.93"
$ws.Range("E49").Value = "  -1.57%  "
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  +0.23%  "
